$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.302.47"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "'2.365.96"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'309.71"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'103.66"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("D7").Value = "'0.511"
$ws.Range("E7").Value = "  -4.61%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "'35.63"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "'53.13"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "'6.92"
$ws.Range("E14").Value = "  -3.90%  "
$ws.Range("D15").Value = "'2.737.82"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").Value = "'15.49"
$ws.Range("E16").Value = "  +3.48%  "
$ws.Range("D17").Value = "'2.365.81"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D19").Value = "'43.293.09"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "'6.33"
$ws.Range("E20").Value = "  +3.44%  "
$ws.Range("E21").Value = "  -5.78%  "
$ws.Range("D22").Value = "'0.0₃0912"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "'68.05"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "'239.70"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "'2.03"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'2.59"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D28").Value = "'25.76"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("D30").Value = "'2.33"
$ws.Range("E30").Value = "  +10.03%  "
$ws.Range("D31").Value = "'36.56"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").Value = "'9.44"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").Value = "'161.48"
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'18.13"
$ws.Range("D37").Value = "'2.49"
$ws.Range("E37").Value = "  +4.15%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'3.07"
$ws.Range("E38").Value = "  -2.89%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'4.66"
$ws.Range("E39").Value = "  +8.66%  "
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("D44").Value = "'2.59"
$ws.Range("E44").Value = "  +12.04%  "
$ws.Range("D45").Value = "'2.039.79"
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("D46").Value = "'19.46"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'0.0288"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "'10.59"
$ws.Range("E48").Value = "  +7.86%  "
$ws.Range("D49").Value = "'3.08"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("D50").Value = "'57.69"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("E51").Value = "  -1.82%  "
